$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fluoxetine_neurons")

# 1) Shrink header row heights
$ws.Rows.Item(2).RowHeight = 73.15
$ws.Rows.Item(3).RowHeight = 18.75

# 2) Hide columns E, G, H (widths stay the same)
$ws.Columns.Item(5).Hidden = $true
$ws.Columns.Item(7).Hidden = $true
$ws.Columns.Item(8).Hidden = $true

# 3) Fix up the label in A21: was the misspelled "mebrane barrier",
#    retype it correctly as "membrane barrier"
$ws.Range("A21").Value = "membrane barrier"

# 4) Delete the now-unneeded detail rows 30:37 (shell110 .. shell145),
#    shifting the summary row up from 40 to 32
$ws.Range("A30:A37").EntireRow.Delete()

# 5) The merged range's last row lost its bottom border when row 37
#    (which used to carry it) was removed - restore it from row 19,
#    which has the same "bottom of group" boxed style
$ws.Range("A19").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 6) Clear the stale SUM formula that used to live in F40 (now F32)
$ws.Range("F32").ClearContents()

# 7) Leave the selection on A21, matching the saved view
$ws.Range("A21").Select()
